$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers (e.g. "44.65") would be
# auto-coerced to a Number by Excel on assignment, which both changes the
# cell type away from the original inlineStr/Text and can introduce binary
# floating point drift (e.g. "86.50" -> 86.5, "0.0591" -> 5.91E-2).
# The source data keeps these as literal text, so force text entry using a
# leading apostrophe (exactly like typing '44.65 into a cell in the UI),
# then reset the cell style back to Normal so no stray "quote prefix" number
# format sticks around on the cell (keeping styles identical to the original).
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

$ws.Range("D2").Value = '28.587.43'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '1.577.05'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.19%  '
Set-TextValue "D5" '213.31'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  +0.18%  '
Set-TextValue "D8" '44.65'
$ws.Range("E8").Value = '  +1.27%  '
Set-TextValue "D9" '24.01'
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("E10").Value = '  -1.49%  '
Set-TextValue "D11" '0.0591'
$ws.Range("E11").Value = '  -1.35%  '
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '1.801.49'
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '1.576.80'
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '28.573.15'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D16" '0.519'
$ws.Range("E16").Value = '  -1.96%  '
$ws.Range("E17").Value = '  -1.70%  '
Set-TextValue "D18" '62.13'
$ws.Range("E18").Value = '  -1.66%  '
Set-TextValue "D19" '230.56'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  -2.18%  '
$ws.Range("E22").Value = '  +0.29%  '
Set-TextValue "D23" '3.89'
$ws.Range("E23").Value = '  -4.34%  '
Set-TextValue "D24" '9.16'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("E25").Value = '  +5.80%  '
Set-TextValue "D26" '151.26'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("E30").Value = '  +0.19%  '
Set-TextValue "D31" '0.0482'
$ws.Range("E31").Value = '  +2.22%  '
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -2.19%  '
$ws.Range("D35").Value = '1.397.29'
$ws.Range("E35").Value = '  -0.11%  '
Set-TextValue "D36" '1.06'
$ws.Range("E36").Value = '  +4.61%  '
$ws.Range("E37").Value = '  -4.09%  '
$ws.Range("E38").Value = '  +0.71%  '
Set-TextValue "D40" '0.0166'
$ws.Range("E40").Value = '  -0.62%  '
Set-TextValue "D41" '0.521'
$ws.Range("E41").Value = '  -3.72%  '
$ws.Range("E42").Value = '  +0.21%  '
Set-TextValue "D43" '0.793'
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("E44").Value = '  -0.25%  '
Set-TextValue "D45" '0.0464'
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("E46").Value = '  -1.92%  '
Set-TextValue "D47" '0.962'
Set-TextValue "D48" '62.95'
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("D49").Value = '1.713.69'
$ws.Range("E49").Value = '  -0.70%  '
Set-TextValue "D50" '86.50'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  -1.16%  '
